$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1971.25
$ws.Range("I40").Value = 1531.1111
$ws.Range("J40").Value = 2143.4783
$ws.Range("K40").Value = 1531.1111
$ws.Range("L40").Value = 2143.4783
$ws.Range("M40").Value = -1356.1111
$ws.Range("N40").Value = -2493.4783
$ws.Range("H107").Value = 941
$ws.Range("I107").Value = 924.75
$ws.Range("J107").Value = 1006
$ws.Range("K107").Value = 924.75
$ws.Range("L107").Value = 1006
$ws.Range("M107").Value = 995.25
$ws.Range("N107").Value = -4846
$ws.Range("H113").Value = 3533.5557
$ws.Range("I113").Value = 3583.75
$ws.Range("J113").Value = 3433.1667
$ws.Range("K113").Value = 3583.75
$ws.Range("L113").Value = 3433.1667
$ws.Range("M113").Value = -329.75
$ws.Range("N113").Value = -9941.1667
$ws.Range("H116").Value = 3895.25
$ws.Range("I116").Value = 4650
$ws.Range("J116").Value = 3517.875
$ws.Range("K116").Value = 4650
$ws.Range("L116").Value = 3517.875
$ws.Range("M116").Value = -1208
$ws.Range("N116").Value = -10401.875
$ws.Range("H132").Value = 8726.25
$ws.Range("I132").Value = 8726.25
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 26178.75
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -23648.75
$ws.Range("N132").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = -84
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents() | Out-Null
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 33334950
$ws.Range("I99").Value = 45456130
$ws.Range("J99").Value = 1705.5
$ws.Range("K99").Value = 45456130
$ws.Range("L99").Value = 1705.5
$ws.Range("M99").Value = -45454632
$ws.Range("N99").Value = -4701.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 4902192.5
$ws.Range("I2").Value = 9804120
$ws.Range("J2").Value = 265.66666
$ws.Range("K2").Value = 58824720
$ws.Range("L2").Value = 1593.99996
$ws.Range("M2").Value = -58824607
$ws.Range("N2").Value = -1819.99996
$ws.Range("H5").Value = 1977.5555
$ws.Range("I5").Value = 2676.3333
$ws.Range("J5").Value = 580
$ws.Range("K5").Value = 8028.999899999999
$ws.Range("L5").Value = 1740
$ws.Range("M5").Value = -7916.999899999999
$ws.Range("N5").Value = -1964
$ws.Range("H12").Value = 46.357143
$ws.Range("I12").Value = 52
$ws.Range("J12").Value = 43.22222
$ws.Range("K12").Value = 156
$ws.Range("L12").Value = 129.66666
$ws.Range("M12").Value = 17
$ws.Range("N12").Value = -475.66666
$ws.Range("H34").Value = 1806.0769
$ws.Range("I34").Value = 1000
$ws.Range("J34").Value = 2746.5
$ws.Range("K34").Value = 3000
$ws.Range("L34").Value = 8239.5
$ws.Range("M34").Value = -2916
$ws.Range("N34").Value = -8407.5
$ws.Range("H38").Value = 680.375
$ws.Range("I38").Value = 97
$ws.Range("K38").Value = 291
$ws.Range("M38").Value = 56
$ws.Range("H39").Value = 5400
$ws.Range("J39").Value = 5400
$ws.Range("L39").Value = 16200
$ws.Range("N39").Value = -16788
$ws.Range("H55").Value = 40894.2
$ws.Range("I55").Value = 250575
$ws.Range("J55").Value = 955
$ws.Range("K55").Value = 751725
$ws.Range("L55").Value = 2865
$ws.Range("M55").Value = -751548
$ws.Range("N55").Value = -3219
$ws.Range("H81").Value = 2333.3333
$ws.Range("I81").Value = 2000
$ws.Range("K81").Value = 6000
$ws.Range("M81").Value = -4877
$ws.Range("H84").Value = 2333.3333
$ws.Range("I84").Value = 2000
$ws.Range("K84").Value = 18000
$ws.Range("M84").Value = -12384
$ws.Range("H131").Value = 8638965
$ws.Range("I131").Value = 71573210
$ws.Range("J131").Value = 932.62744
$ws.Range("K131").Value = 214719630
$ws.Range("L131").Value = 2797.88232
$ws.Range("M131").Value = -214714590
$ws.Range("N131").Value = -12877.88232
$ws.Range("H135").Value = 1977.5555
$ws.Range("I135").Value = 2676.3333
$ws.Range("J135").Value = 580
$ws.Range("K135").Value = 24086.9997
$ws.Range("L135").Value = 5220
$ws.Range("M135").Value = -21551.9997
$ws.Range("N135").Value = -10290

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5366.5
$ws.Range("I70").Value = 5649.75
$ws.Range("K70").Value = 5649.75
$ws.Range("M70").Value = -5379.75
$ws.Range("H73").Value = 5366.5
$ws.Range("I73").Value = 5649.75
$ws.Range("K73").Value = 5649.75
$ws.Range("M73").Value = -4713.75
$ws.Range("H113").Value = 6492.1
$ws.Range("I113").Value = 9802.200000000001
$ws.Range("J113").Value = 3182
$ws.Range("K113").Value = 9802.200000000001
$ws.Range("L113").Value = 3182
$ws.Range("M113").Value = -7632.200000000001
$ws.Range("N113").Value = -7522

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 450
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -1090
$ws.Range("H27").Value = 450
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -714
$ws.Range("H61").Value = 8295.799999999999
$ws.Range("I61").Value = 10902.909
$ws.Range("J61").Value = 1126.25
$ws.Range("K61").Value = 10902.909
$ws.Range("L61").Value = 1126.25
$ws.Range("M61").Value = -10700.909
$ws.Range("N61").Value = -1530.25
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents() | Out-Null
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents() | Out-Null
$ws.Range("H113").Value = 8295.799999999999
$ws.Range("I113").Value = 10902.909
$ws.Range("J113").Value = 1126.25
$ws.Range("K113").Value = 10902.909
$ws.Range("L113").Value = 1126.25
$ws.Range("M113").Value = -8732.909
$ws.Range("N113").Value = -5466.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 858.871
$ws.Range("I107").Value = 910.95654
$ws.Range("J107").Value = 709.125
$ws.Range("K107").Value = 2732.86962
$ws.Range("L107").Value = 2127.375
$ws.Range("M107").Value = -812.8696199999999
$ws.Range("N107").Value = -5967.375
$ws.Range("H113").Value = 572.7857
$ws.Range("I113").Value = 351.72726
$ws.Range("J113").Value = 1383.3334
$ws.Range("K113").Value = 1055.18178
$ws.Range("L113").Value = 4150.0002
$ws.Range("M113").Value = 1114.81822
$ws.Range("N113").Value = -8490.0002
